$d = $word.ActiveDocument

$replacements = @(
    @{Old="2025-07-08 Tuesday"; New="2025-07-09 Wednesday"},
    @{Old="359×3="; New="935×5="},
    @{Old="239×7="; New="157×8="},
    @{Old="647×2="; New="769×7="},
    @{Old="716×6="; New="177×6="},
    @{Old="299×3="; New="611×2="},
    @{Old="338×3="; New="498×9="},
    @{Old="808×2="; New="458×6="},
    @{Old="731×8="; New="892×8="},
    @{Old="757×3="; New="554×8="},
    @{Old="962×6="; New="222×9="},
    @{Old="813×8="; New="852×5="},
    @{Old="401×4="; New="343×3="},
    @{Old="822×4="; New="773×4="},
    @{Old="239×6="; New="699×3="},
    @{Old="720×8="; New="289×6="},
    @{Old="833×9="; New="603×8="},
    @{Old="762×2="; New="727×3="},
    @{Old="450×7="; New="952×3="},
    @{Old="473×7="; New="636×2="},
    @{Old="852×4="; New="414×8="},
    @{Old="480×7="; New="627×6="},
    @{Old="766×3="; New="384×3="},
    @{Old="733×7="; New="892×2="},
    @{Old="586×2="; New="602×6="},
    @{Old="885×6="; New="288×2="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
